$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.337.78"
$ws.Range("E2").Value = "  -0.86%  "
$ws.Range("D3").Value = "'1.881.67"
$ws.Range("E3").Value = "  -0.65%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("D5").Value = "'237.51"
$ws.Range("E5").Value = "  -1.03%  "
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("D7").Value = "'0.4821"
$ws.Range("E7").Value = "  -2.20%  "
$ws.Range("D8").Value = "'0.2860"
$ws.Range("E8").Value = "  -3.05%  "
$ws.Range("D9").Value = "'0.06586"
$ws.Range("E9").Value = "  -2.07%  "
$ws.Range("D10").Value = "'1.885.01"
$ws.Range("E10").Value = "  -0.29%  "
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").Value = "'0.07405"
$ws.Range("E11").Value = "  +0.58%  "
$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").Value = "'16.66"
$ws.Range("E12").Value = "  -3.41%  "
$ws.Range("D13").Value = "'5.213"
$ws.Range("E13").Value = "  +0.90%  "
$ws.Range("D14").Value = "'88.11"
$ws.Range("E14").Value = "  -0.30%  "
$ws.Range("D15").Value = "'0.6558"
$ws.Range("E15").Value = "  -2.46%  "
$ws.Range("D16").Value = "'30.328.90"
$ws.Range("E16").Value = "  -0.69%  "
$ws.Range("D17").Value = "'13.43"
$ws.Range("E17").Value = "  -0.13%  "
$ws.Range("B18").Value = "Dai"
$ws.Range("C18").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D18").Value = "'1.003"
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "'0.000007713"
$ws.Range("E19").Value = "  -2.17%  "
$ws.Range("D20").Value = "'2.135.63"
$ws.Range("E20").Value = "  -1.92%  "
$ws.Range("D21").Value = "'5.358"
$ws.Range("E21").Value = "  +0.52%  "
$ws.Range("D22").Value = "'1.004"
$ws.Range("E22").Value = "  +0.28%  "
$ws.Range("D23").Value = "'218.95"
$ws.Range("E23").Value = "  +14.46%  "
$ws.Range("D24").Value = "'6.160"
$ws.Range("E24").Value = "  -1.04%  "
$ws.Range("D25").Value = "'9.387"
$ws.Range("E25").Value = "  -1.78%  "
$ws.Range("D26").Value = "'165.34"
$ws.Range("E26").Value = "  +2.66%  "
$ws.Range("D27").Value = "'18.88"
$ws.Range("E27").Value = "  +2.36%  "
$ws.Range("D28").Value = "'1.929"
$ws.Range("E28").Value = "  -0.86%  "
$ws.Range("D29").Value = "'1.459"
$ws.Range("E29").Value = "  -0.98%  "
$ws.Range("D30").Value = "'4.293"
$ws.Range("E30").Value = "  -3.57%  "
$ws.Range("D31").Value = "'0.09244"
$ws.Range("E31").Value = "  +0.55%  "
$ws.Range("D32").Value = "'4.009"
$ws.Range("E32").Value = "  -4.85%  "
$ws.Range("D33").Value = "'0.05065"
$ws.Range("E33").Value = "  -3.57%  "
$ws.Range("D34").Value = "'0.7515"
$ws.Range("E34").Value = "  +0.99%  "
$ws.Range("D35").Value = "'1.144"
$ws.Range("E35").Value = "  +3.50%  "
$ws.Range("E36").Value = "  -0.62%  "
$ws.Range("D37").Value = "'0.01846"
$ws.Range("E37").Value = "  +0.59%  "
$ws.Range("D38").Value = "'2.639"
$ws.Range("E38").Value = "  -2.00%  "
$ws.Range("D39").Value = "'0.9145"
$ws.Range("E39").Value = "  -0.90%  "
$ws.Range("D40").Value = "'2.075"
$ws.Range("E40").Value = "  +0.32%  "
$ws.Range("D41").Value = "'107.35"
$ws.Range("E41").Value = "  +1.07%  "
$ws.Range("D42").Value = "'5.946"
$ws.Range("E42").Value = "  -0.15%  "
$ws.Range("D43").Value = "'0.4307"
$ws.Range("E43").Value = "  -2.98%  "
$ws.Range("D44").Value = "'1.003"
$ws.Range("E44").Value = "  +0.94%  "
$ws.Range("D45").Value = "'7.575"
$ws.Range("E45").Value = "  -0.29%  "
$ws.Range("D46").Value = "'1.581"
$ws.Range("E46").Value = "  +10.97%  "
$ws.Range("D47").Value = "'0.1320"
$ws.Range("E47").Value = "  -4.96%  "
$ws.Range("D48").Value = "'64.98"
$ws.Range("E48").Value = "  -9.77%  "
$ws.Range("D49").Value = "'8.904"
$ws.Range("E49").Value = "  -1.37%  "
$ws.Range("D50").Value = "'34.58"
$ws.Range("E50").Value = "  -1.58%  "
$ws.Range("D51").Value = "'0.05725"
$ws.Range("E51").Value = "  -1.86%  "
